$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows at 377-378; this shifts the existing rows 377-428 down
# to 379-430, which matches the new dimension A1:T430.
$ws.Rows("377:378").Insert()

# New row 377 - "Especial" quality entry for the new reporting week.
$ws.Range("A377").Value = 3
$ws.Range("B377").Value = "Femacal de La Calera"
$ws.Range("C377").Value = "Coquimbo"
$ws.Range("D377").Value = 45212
$ws.Range("E377").Value = 5
$ws.Range("F377").Value = "Fruta"
$ws.Range("G377").Value = 100107
$ws.Range("H377").Value = "Otros"
$ws.Range("I377").Value = 100107002
$ws.Range("J377").Value = "Chirimoya"
$ws.Range("K377").Value = "Cultivar IV Región"
$ws.Range("L377").Value = "Especial"
$ws.Range("M377").Value = 92
$ws.Range("N377").Value = 28000
$ws.Range("O377").Value = 30000
$ws.Range("P377").Value = 28978
$ws.Range("Q377").Value = "$/bandeja 10 kilos"
$ws.Range("R377").Value = "Provincia del Elquí"
$ws.Range("S377").Value = 2898
$ws.Range("T377").Value = 10

# New row 378 - "Primera" quality entry for the new reporting week.
$ws.Range("A378").Value = 3
$ws.Range("B378").Value = "Femacal de La Calera"
$ws.Range("C378").Value = "Coquimbo"
$ws.Range("D378").Value = 45212
$ws.Range("E378").Value = 5
$ws.Range("F378").Value = "Fruta"
$ws.Range("G378").Value = 100107
$ws.Range("H378").Value = "Otros"
$ws.Range("I378").Value = 100107002
$ws.Range("J378").Value = "Chirimoya"
$ws.Range("K378").Value = "Cultivar IV Región"
$ws.Range("L378").Value = "Primera"
$ws.Range("M378").Value = 90
$ws.Range("N378").Value = 25000
$ws.Range("O378").Value = 26000
$ws.Range("P378").Value = 25444
$ws.Range("Q378").Value = "$/bandeja 10 kilos"
$ws.Range("R378").Value = "Provincia del Elquí"
$ws.Range("S378").Value = 2544
$ws.Range("T378").Value = 10
